# Refresh the "cryptos" price/volume snapshot (GitHub Actions scheduled update).
# Price (col D) and 1h volume-change (col E) are plain text cells, not numbers
# (values like "26.985.27" use dots as thousands separators, and percentages
# keep their padding spaces), so every write below targets Range.Value with a
# string. For prices that would otherwise auto-parse as a real number (e.g.
# "1.005", "91.20", "0.08900") a leading apostrophe forces Excel to keep the
# literal text (preserving trailing zeros / exact formatting); the style is
# reset back to "Normal" right after so the cell's quote-prefix indicator
# doesn't leave a lingering explicit style on an originally unstyled cell.
# Rows 46/47 additionally swap which coin (EnergySwap / PaxDollar) occupies
# which rank row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.968.53"
$ws.Range("E2").Value = "  -0.55%  "

$ws.Range("D3").Value = "1.824.03"
$ws.Range("E3").Value = "  -0.06%  "

$ws.Range("D4").Value = "'1.005"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.44%  "

$ws.Range("D5").Value = "'311.32"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.01%  "

$ws.Range("D6").Value = "'1.004"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.44%  "

$ws.Range("D7").Value = "'0.4629"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.27%  "

$ws.Range("D8").Value = "'0.3707"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.89%  "

$ws.Range("D9").Value = "'0.07345"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.67%  "

$ws.Range("D10").Value = "'0.8742"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.49%  "

$ws.Range("D11").Value = "'0.07911"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +3.67%  "

$ws.Range("D12").Value = "'19.76"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.08%  "

$ws.Range("D13").Value = "1.864.59"
$ws.Range("E13").Value = "  +1.85%  "

$ws.Range("D14").Value = "'5.333"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.22%  "

$ws.Range("D15").Value = "'6.546"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.99%  "

$ws.Range("D16").Value = "'91.20"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.64%  "

$ws.Range("D17").Value = "'1.006"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.24%  "

$ws.Range("D18").Value = "'0.000008851"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.29%  "

$ws.Range("E19").Value = "  -0.31%  "

$ws.Range("D20").Value = "'14.77"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.86%  "

$ws.Range("D21").Value = "27.006.28"
$ws.Range("E21").Value = "  -0.88%  "

$ws.Range("D22").Value = "'5.096"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.80%  "

$ws.Range("D23").Value = "'10.53"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.33%  "

$ws.Range("D24").Value = "2.086.98"
$ws.Range("E24").Value = "  +0.17%  "

$ws.Range("D25").Value = "'153.31"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.03%  "

$ws.Range("D26").Value = "'1.857"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.34%  "

$ws.Range("D27").Value = "'18.37"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.50%  "

$ws.Range("D28").Value = "'2.037"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -3.46%  "

$ws.Range("D29").Value = "'5.127"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.76%  "

$ws.Range("D30").Value = "'115.59"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.49%  "

$ws.Range("D31").Value = "'0.08900"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.17%  "

$ws.Range("D32").Value = "'2.958"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.03%  "

$ws.Range("D33").Value = "'0.7277"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.07%  "

$ws.Range("D34").Value = "'4.435"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.33%  "

$ws.Range("E35").Value = "  -1.00%  "

$ws.Range("E36").Value = "  -3.09%  "

$ws.Range("E37").Value = "  +0.34%  "

$ws.Range("D39").Value = "'0.05216"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.95%  "

$ws.Range("D40").Value = "'2.948"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.50%  "

$ws.Range("D41").Value = "'7.101"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.46%  "

$ws.Range("D42").Value = "'0.5149"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.39%  "

$ws.Range("D43").Value = "'0.1620"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.85%  "

$ws.Range("D44").Value = "'8.168"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.26%  "

$ws.Range("D45").Value = "'0.4830"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.97%  "

$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").Value = "'10.22"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.81%  "

$ws.Range("B47").Value = "PaxDollar"
$ws.Range("C47").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D47").Value = "'1.004"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.49%  "

$ws.Range("D48").Value = "'102.57"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.36%  "

$ws.Range("D49").Value = "'1.632"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.36%  "

$ws.Range("D50").Value = "'0.06193"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.82%  "

$ws.Range("D51").Value = "'64.90"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.65%  "
